$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ============ Part 1: swap paired rows (columns B:AC) ============
# rows 38 <-> 39
$ws.Cells.Item(38,2).Value = 6782522
$ws.Cells.Item(38,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(38,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(38,5).Value = 45171.75
$ws.Cells.Item(38,6).Value = 'Municipal Perez Zeledon'
$ws.Cells.Item(38,7).Value = 'Sporting San Jose'
$ws.Cells.Item(38,8).Value = 1
$ws.Cells.Item(38,9).Value = 2
$ws.Cells.Item(38,10).Value = 'A'
$ws.Cells.Item(38,11).Value = 2.5
$ws.Cells.Item(38,12).Value = 3.5
$ws.Cells.Item(38,13).Value = 2.5
$ws.Cells.Item(38,14).Value = 2.2
$ws.Cells.Item(38,15).Value = 3.5
$ws.Cells.Item(38,16).Value = 2.9
$ws.Cells.Item(38,17).Value = -0.25
$ws.Cells.Item(38,18).Value = 1.9
$ws.Cells.Item(38,19).Value = 1.9
$ws.Cells.Item(38,20).Value = 2.5
$ws.Cells.Item(38,21).Value = 1.9
$ws.Cells.Item(38,22).Value = 1.9
$ws.Cells.Item(38,23).Value = -1
$ws.Cells.Item(38,24).Value = -1
$ws.Cells.Item(38,25).Value = 1.9
$ws.Cells.Item(38,26).Value = -1
$ws.Cells.Item(38,27).Value = 0.8999999999999999
$ws.Cells.Item(38,28).Value = 0.8999999999999999
$ws.Cells.Item(38,29).Value = -1
$ws.Cells.Item(39,2).Value = 6781354
$ws.Cells.Item(39,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(39,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(39,5).Value = 45171.75
$ws.Cells.Item(39,6).Value = 'Puntarenas'
$ws.Cells.Item(39,7).Value = 'AD San Carlos'
$ws.Cells.Item(39,8).Value = 1
$ws.Cells.Item(39,9).Value = 0
$ws.Cells.Item(39,10).Value = 'H'
$ws.Cells.Item(39,11).Value = 2.4
$ws.Cells.Item(39,12).Value = 3.2
$ws.Cells.Item(39,13).Value = 2.8
$ws.Cells.Item(39,14).Value = 2.3
$ws.Cells.Item(39,15).Value = 3.2
$ws.Cells.Item(39,16).Value = 3
$ws.Cells.Item(39,17).Value = -0.25
$ws.Cells.Item(39,18).Value = 2
$ws.Cells.Item(39,19).Value = 1.8
$ws.Cells.Item(39,20).Value = 2.25
$ws.Cells.Item(39,21).Value = 1.9
$ws.Cells.Item(39,22).Value = 1.9
$ws.Cells.Item(39,23).Value = 1.3
$ws.Cells.Item(39,24).Value = -1
$ws.Cells.Item(39,25).Value = -1
$ws.Cells.Item(39,26).Value = 1
$ws.Cells.Item(39,27).Value = -1
$ws.Cells.Item(39,28).Value = -1
$ws.Cells.Item(39,29).Value = 0.8999999999999999

# rows 91 <-> 92
$ws.Cells.Item(91,2).Value = 6782568
$ws.Cells.Item(91,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(91,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(91,5).Value = 45220.83333333334
$ws.Cells.Item(91,6).Value = 'Sporting San Jose'
$ws.Cells.Item(91,7).Value = 'AD Guanacasteca'
$ws.Cells.Item(91,8).Value = 1
$ws.Cells.Item(91,9).Value = 1
$ws.Cells.Item(91,10).Value = 'D'
$ws.Cells.Item(91,11).Value = 1.909
$ws.Cells.Item(91,12).Value = 3.6
$ws.Cells.Item(91,13).Value = 3.3
$ws.Cells.Item(91,14).Value = 2
$ws.Cells.Item(91,15).Value = 3.6
$ws.Cells.Item(91,16).Value = 3.1
$ws.Cells.Item(91,17).Value = -0.5
$ws.Cells.Item(91,18).Value = 2
$ws.Cells.Item(91,19).Value = 1.8
$ws.Cells.Item(91,20).Value = 2.5
$ws.Cells.Item(91,21).Value = 1.825
$ws.Cells.Item(91,22).Value = 1.975
$ws.Cells.Item(91,23).Value = -1
$ws.Cells.Item(91,24).Value = 2.6
$ws.Cells.Item(91,25).Value = -1
$ws.Cells.Item(91,26).Value = -1
$ws.Cells.Item(91,27).Value = 0.8
$ws.Cells.Item(91,28).Value = -1
$ws.Cells.Item(91,29).Value = 0.9750000000000001
$ws.Cells.Item(92,2).Value = 6782566
$ws.Cells.Item(92,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(92,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(92,5).Value = 45220.83333333334
$ws.Cells.Item(92,6).Value = 'Cartagines'
$ws.Cells.Item(92,7).Value = 'Deportivo Saprissa'
$ws.Cells.Item(92,8).Value = 0
$ws.Cells.Item(92,9).Value = 4
$ws.Cells.Item(92,10).Value = 'A'
$ws.Cells.Item(92,11).Value = 3.2
$ws.Cells.Item(92,12).Value = 3.4
$ws.Cells.Item(92,13).Value = 2
$ws.Cells.Item(92,14).Value = 2.9
$ws.Cells.Item(92,15).Value = 3.5
$ws.Cells.Item(92,16).Value = 2.15
$ws.Cells.Item(92,17).Value = 0.25
$ws.Cells.Item(92,18).Value = 1.875
$ws.Cells.Item(92,19).Value = 1.925
$ws.Cells.Item(92,20).Value = 3
$ws.Cells.Item(92,21).Value = 1.975
$ws.Cells.Item(92,22).Value = 1.825
$ws.Cells.Item(92,23).Value = -1
$ws.Cells.Item(92,24).Value = -1
$ws.Cells.Item(92,25).Value = 1.15
$ws.Cells.Item(92,26).Value = -1
$ws.Cells.Item(92,27).Value = 0.925
$ws.Cells.Item(92,28).Value = 0.9750000000000001
$ws.Cells.Item(92,29).Value = -1

# rows 95 <-> 96
$ws.Cells.Item(95,2).Value = 6782567
$ws.Cells.Item(95,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(95,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(95,5).Value = 45221.79166666666
$ws.Cells.Item(95,6).Value = 'AD Grecia'
$ws.Cells.Item(95,7).Value = 'Municipal Liberia'
$ws.Cells.Item(95,8).Value = 2
$ws.Cells.Item(95,9).Value = 3
$ws.Cells.Item(95,10).Value = 'A'
$ws.Cells.Item(95,11).Value = 2.875
$ws.Cells.Item(95,12).Value = 3.5
$ws.Cells.Item(95,13).Value = 2.15
$ws.Cells.Item(95,14).Value = 2.3
$ws.Cells.Item(95,15).Value = 3.5
$ws.Cells.Item(95,16).Value = 2.6
$ws.Cells.Item(95,17).Value = 0
$ws.Cells.Item(95,18).Value = 1.8
$ws.Cells.Item(95,19).Value = 2
$ws.Cells.Item(95,20).Value = 2.75
$ws.Cells.Item(95,21).Value = 1.8
$ws.Cells.Item(95,22).Value = 2
$ws.Cells.Item(95,23).Value = -1
$ws.Cells.Item(95,24).Value = -1
$ws.Cells.Item(95,25).Value = 1.6
$ws.Cells.Item(95,26).Value = -1
$ws.Cells.Item(95,27).Value = 1
$ws.Cells.Item(95,28).Value = 0.8
$ws.Cells.Item(95,29).Value = -1
$ws.Cells.Item(96,2).Value = 6782565
$ws.Cells.Item(96,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(96,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(96,5).Value = 45221.79166666666
$ws.Cells.Item(96,6).Value = 'Santos de Gupiles'
$ws.Cells.Item(96,7).Value = 'Municipal Perez Zeledon'
$ws.Cells.Item(96,8).Value = 2
$ws.Cells.Item(96,9).Value = 0
$ws.Cells.Item(96,10).Value = 'H'
$ws.Cells.Item(96,11).Value = 1.833
$ws.Cells.Item(96,12).Value = 3.4
$ws.Cells.Item(96,13).Value = 3.75
$ws.Cells.Item(96,14).Value = 1.833
$ws.Cells.Item(96,15).Value = 3.5
$ws.Cells.Item(96,16).Value = 3.75
$ws.Cells.Item(96,17).Value = -0.5
$ws.Cells.Item(96,18).Value = 1.875
$ws.Cells.Item(96,19).Value = 1.925
$ws.Cells.Item(96,20).Value = 2.75
$ws.Cells.Item(96,21).Value = 2
$ws.Cells.Item(96,22).Value = 1.8
$ws.Cells.Item(96,23).Value = 0.833
$ws.Cells.Item(96,24).Value = -1
$ws.Cells.Item(96,25).Value = -1
$ws.Cells.Item(96,26).Value = 0.875
$ws.Cells.Item(96,27).Value = -1
$ws.Cells.Item(96,28).Value = -1
$ws.Cells.Item(96,29).Value = 0.8

# rows 110 <-> 111
$ws.Cells.Item(110,2).Value = 6782579
$ws.Cells.Item(110,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(110,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(110,5).Value = 45238.875
$ws.Cells.Item(110,6).Value = 'Santos de Gupiles'
$ws.Cells.Item(110,7).Value = 'AD San Carlos'
$ws.Cells.Item(110,8).Value = 0
$ws.Cells.Item(110,9).Value = 2
$ws.Cells.Item(110,10).Value = 'A'
$ws.Cells.Item(110,11).Value = 2.4
$ws.Cells.Item(110,12).Value = 3.3
$ws.Cells.Item(110,13).Value = 2.7
$ws.Cells.Item(110,14).Value = 2.375
$ws.Cells.Item(110,15).Value = 3.4
$ws.Cells.Item(110,16).Value = 2.8
$ws.Cells.Item(110,17).Value = -0.25
$ws.Cells.Item(110,18).Value = 2
$ws.Cells.Item(110,19).Value = 1.8
$ws.Cells.Item(110,20).Value = 2.5
$ws.Cells.Item(110,21).Value = 1.875
$ws.Cells.Item(110,22).Value = 1.925
$ws.Cells.Item(110,23).Value = -1
$ws.Cells.Item(110,24).Value = -1
$ws.Cells.Item(110,25).Value = 1.8
$ws.Cells.Item(110,26).Value = -1
$ws.Cells.Item(110,27).Value = 0.8
$ws.Cells.Item(110,28).Value = -1
$ws.Cells.Item(110,29).Value = 0.925
$ws.Cells.Item(111,2).Value = 6782581
$ws.Cells.Item(111,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(111,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(111,5).Value = 45238.875
$ws.Cells.Item(111,6).Value = 'Alajuelense'
$ws.Cells.Item(111,7).Value = 'AD Grecia'
$ws.Cells.Item(111,8).Value = 2
$ws.Cells.Item(111,9).Value = 0
$ws.Cells.Item(111,10).Value = 'H'
$ws.Cells.Item(111,11).Value = 1.181
$ws.Cells.Item(111,12).Value = 6.5
$ws.Cells.Item(111,13).Value = 11
$ws.Cells.Item(111,14).Value = 1.25
$ws.Cells.Item(111,15).Value = 5
$ws.Cells.Item(111,16).Value = 9
$ws.Cells.Item(111,17).Value = -1.75
$ws.Cells.Item(111,18).Value = 1.975
$ws.Cells.Item(111,19).Value = 1.825
$ws.Cells.Item(111,20).Value = 3.25
$ws.Cells.Item(111,21).Value = 2
$ws.Cells.Item(111,22).Value = 1.8
$ws.Cells.Item(111,23).Value = 0.25
$ws.Cells.Item(111,24).Value = -1
$ws.Cells.Item(111,25).Value = -1
$ws.Cells.Item(111,26).Value = 0.4875
$ws.Cells.Item(111,27).Value = -0.5
$ws.Cells.Item(111,28).Value = -1
$ws.Cells.Item(111,29).Value = 0.8

# rows 130 <-> 131
$ws.Cells.Item(130,2).Value = 6782596
$ws.Cells.Item(130,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(130,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(130,5).Value = 45255.95833333334
$ws.Cells.Item(130,6).Value = 'Alajuelense'
$ws.Cells.Item(130,7).Value = 'AD Guanacasteca'
$ws.Cells.Item(130,8).Value = 3
$ws.Cells.Item(130,9).Value = 4
$ws.Cells.Item(130,10).Value = 'A'
$ws.Cells.Item(130,11).Value = 1.363
$ws.Cells.Item(130,12).Value = 4.75
$ws.Cells.Item(130,13).Value = 8
$ws.Cells.Item(130,14).Value = 1.444
$ws.Cells.Item(130,15).Value = 4.333
$ws.Cells.Item(130,16).Value = 7
$ws.Cells.Item(130,17).Value = -1.25
$ws.Cells.Item(130,18).Value = 1.975
$ws.Cells.Item(130,19).Value = 1.825
$ws.Cells.Item(130,20).Value = 2.75
$ws.Cells.Item(130,21).Value = 1.775
$ws.Cells.Item(130,22).Value = 2.025
$ws.Cells.Item(130,23).Value = -1
$ws.Cells.Item(130,24).Value = -1
$ws.Cells.Item(130,25).Value = 6
$ws.Cells.Item(130,26).Value = -1
$ws.Cells.Item(130,27).Value = 0.825
$ws.Cells.Item(130,28).Value = 0.7749999999999999
$ws.Cells.Item(130,29).Value = -1
$ws.Cells.Item(131,2).Value = 6782598
$ws.Cells.Item(131,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(131,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(131,5).Value = 45255.95833333334
$ws.Cells.Item(131,6).Value = 'Municipal Perez Zeledon'
$ws.Cells.Item(131,7).Value = 'Cartagines'
$ws.Cells.Item(131,8).Value = 1
$ws.Cells.Item(131,9).Value = 0
$ws.Cells.Item(131,10).Value = 'H'
$ws.Cells.Item(131,11).Value = 4.5
$ws.Cells.Item(131,12).Value = 3.75
$ws.Cells.Item(131,13).Value = 1.615
$ws.Cells.Item(131,14).Value = 3.4
$ws.Cells.Item(131,15).Value = 3.4
$ws.Cells.Item(131,16).Value = 1.85
$ws.Cells.Item(131,17).Value = 0.5
$ws.Cells.Item(131,18).Value = 1.8
$ws.Cells.Item(131,19).Value = 2
$ws.Cells.Item(131,20).Value = 2.75
$ws.Cells.Item(131,21).Value = 1.9
$ws.Cells.Item(131,22).Value = 1.9
$ws.Cells.Item(131,23).Value = 2.4
$ws.Cells.Item(131,24).Value = -1
$ws.Cells.Item(131,25).Value = -1
$ws.Cells.Item(131,26).Value = 0.8
$ws.Cells.Item(131,27).Value = -1
$ws.Cells.Item(131,28).Value = -1
$ws.Cells.Item(131,29).Value = 0.8999999999999999

# rows 192 <-> 193
$ws.Cells.Item(192,2).Value = 7623916
$ws.Cells.Item(192,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(192,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(192,5).Value = 45347.75
$ws.Cells.Item(192,6).Value = 'Santos de Gupiles'
$ws.Cells.Item(192,7).Value = 'AD Grecia'
$ws.Cells.Item(192,8).Value = 0
$ws.Cells.Item(192,9).Value = 2
$ws.Cells.Item(192,10).Value = 'A'
$ws.Cells.Item(192,11).Value = 2.05
$ws.Cells.Item(192,12).Value = 3.3
$ws.Cells.Item(192,13).Value = 3.2
$ws.Cells.Item(192,14).Value = 1.909
$ws.Cells.Item(192,15).Value = 3.4
$ws.Cells.Item(192,16).Value = 3.6
$ws.Cells.Item(192,17).Value = -0.5
$ws.Cells.Item(192,18).Value = 1.95
$ws.Cells.Item(192,19).Value = 1.85
$ws.Cells.Item(192,20).Value = 2.5
$ws.Cells.Item(192,21).Value = 1.85
$ws.Cells.Item(192,22).Value = 1.95
$ws.Cells.Item(192,23).Value = -1
$ws.Cells.Item(192,24).Value = -1
$ws.Cells.Item(192,25).Value = 2.6
$ws.Cells.Item(192,26).Value = -1
$ws.Cells.Item(192,27).Value = 0.8500000000000001
$ws.Cells.Item(192,28).Value = -1
$ws.Cells.Item(192,29).Value = 0.95
$ws.Cells.Item(193,2).Value = 7623919
$ws.Cells.Item(193,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(193,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(193,5).Value = 45347.75
$ws.Cells.Item(193,6).Value = 'Municipal Liberia'
$ws.Cells.Item(193,7).Value = 'Sporting San Jose'
$ws.Cells.Item(193,8).Value = 2
$ws.Cells.Item(193,9).Value = 0
$ws.Cells.Item(193,10).Value = 'H'
$ws.Cells.Item(193,11).Value = 1.75
$ws.Cells.Item(193,12).Value = 3.6
$ws.Cells.Item(193,13).Value = 3.8
$ws.Cells.Item(193,14).Value = 1.8
$ws.Cells.Item(193,15).Value = 3.6
$ws.Cells.Item(193,16).Value = 3.6
$ws.Cells.Item(193,17).Value = -0.5
$ws.Cells.Item(193,18).Value = 1.9
$ws.Cells.Item(193,19).Value = 1.9
$ws.Cells.Item(193,20).Value = 2.75
$ws.Cells.Item(193,21).Value = 2
$ws.Cells.Item(193,22).Value = 1.8
$ws.Cells.Item(193,23).Value = 0.8
$ws.Cells.Item(193,24).Value = -1
$ws.Cells.Item(193,25).Value = -1
$ws.Cells.Item(193,26).Value = 0.8999999999999999
$ws.Cells.Item(193,27).Value = -1
$ws.Cells.Item(193,28).Value = -1
$ws.Cells.Item(193,29).Value = 0.8
# ============ Part 2: append new rows 199-204 ============
# Row 199
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(199,1).PasteSpecial(-4122)
$ws.Cells.Item(199,1).Value = 197
$ws.Cells.Item(199,2).Value = 7623920
$ws.Cells.Item(199,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(199,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(199,5).PasteSpecial(-4122)
$ws.Cells.Item(199,5).Value = 45352.95833333334
$ws.Cells.Item(199,6).Value = 'Alajuelense'
$ws.Cells.Item(199,7).Value = 'Santos de Gupiles'
$ws.Cells.Item(199,8).Value = 4
$ws.Cells.Item(199,9).Value = 0
$ws.Cells.Item(199,10).Value = 'H'
$ws.Cells.Item(199,11).Value = 1.25
$ws.Cells.Item(199,12).Value = 5.25
$ws.Cells.Item(199,13).Value = 11
$ws.Cells.Item(199,14).Value = 1.166
$ws.Cells.Item(199,15).Value = 7
$ws.Cells.Item(199,16).Value = 17
$ws.Cells.Item(199,17).Value = -2
$ws.Cells.Item(199,18).Value = 1.8
$ws.Cells.Item(199,19).Value = 2
$ws.Cells.Item(199,20).Value = 3.25
$ws.Cells.Item(199,21).Value = 1.95
$ws.Cells.Item(199,22).Value = 1.85
$ws.Cells.Item(199,23).Value = 0.1659999999999999
$ws.Cells.Item(199,24).Value = -1
$ws.Cells.Item(199,25).Value = -1
$ws.Cells.Item(199,26).Value = 0.8
$ws.Cells.Item(199,27).Value = -1
$ws.Cells.Item(199,28).Value = 0.95
$ws.Cells.Item(199,29).Value = -1

# Row 200 (no H/I/J, no AB/AC - match not yet played)
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(200,1).PasteSpecial(-4122)
$ws.Cells.Item(200,1).Value = 198
$ws.Cells.Item(200,2).Value = 7624967
$ws.Cells.Item(200,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(200,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(200,5).PasteSpecial(-4122)
$ws.Cells.Item(200,5).Value = 45353.75
$ws.Cells.Item(200,6).Value = 'Puntarenas'
$ws.Cells.Item(200,7).Value = 'Herediano'
$ws.Cells.Item(200,11).Value = 3.75
$ws.Cells.Item(200,12).Value = 3.4
$ws.Cells.Item(200,13).Value = 1.8
$ws.Cells.Item(200,14).Value = 3.8
$ws.Cells.Item(200,15).Value = 3.4
$ws.Cells.Item(200,16).Value = 1.8
$ws.Cells.Item(200,17).Value = 0.5
$ws.Cells.Item(200,18).Value = 1.9
$ws.Cells.Item(200,19).Value = 1.9
$ws.Cells.Item(200,20).Value = 2.5
$ws.Cells.Item(200,21).Value = 1.95
$ws.Cells.Item(200,22).Value = 1.85
$ws.Cells.Item(200,23).Value = 0
$ws.Cells.Item(200,24).Value = 0
$ws.Cells.Item(200,25).Value = 0
$ws.Cells.Item(200,26).Value = 0
$ws.Cells.Item(200,27).Value = 0

# Row 201
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(201,1).PasteSpecial(-4122)
$ws.Cells.Item(201,1).Value = 199
$ws.Cells.Item(201,2).Value = 7623921
$ws.Cells.Item(201,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(201,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(201,5).PasteSpecial(-4122)
$ws.Cells.Item(201,5).Value = 45353.75
$ws.Cells.Item(201,6).Value = 'AD Grecia'
$ws.Cells.Item(201,7).Value = 'Municipal Liberia'
$ws.Cells.Item(201,11).Value = 2.75
$ws.Cells.Item(201,12).Value = 3.25
$ws.Cells.Item(201,13).Value = 2.3
$ws.Cells.Item(201,14).Value = 2.875
$ws.Cells.Item(201,15).Value = 3.25
$ws.Cells.Item(201,16).Value = 2.25
$ws.Cells.Item(201,17).Value = 0.25
$ws.Cells.Item(201,18).Value = 1.775
$ws.Cells.Item(201,19).Value = 2.025
$ws.Cells.Item(201,20).Value = 2.5
$ws.Cells.Item(201,21).Value = 1.825
$ws.Cells.Item(201,22).Value = 1.975
$ws.Cells.Item(201,23).Value = 0
$ws.Cells.Item(201,24).Value = 0
$ws.Cells.Item(201,25).Value = 0
$ws.Cells.Item(201,26).Value = 0
$ws.Cells.Item(201,27).Value = 0

# Row 202
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(202,1).PasteSpecial(-4122)
$ws.Cells.Item(202,1).Value = 200
$ws.Cells.Item(202,2).Value = 7623922
$ws.Cells.Item(202,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(202,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(202,5).PasteSpecial(-4122)
$ws.Cells.Item(202,5).Value = 45353.95833333334
$ws.Cells.Item(202,6).Value = 'Municipal Perez Zeledon'
$ws.Cells.Item(202,7).Value = 'AD Guanacasteca'
$ws.Cells.Item(202,11).Value = 2.375
$ws.Cells.Item(202,12).Value = 3.4
$ws.Cells.Item(202,13).Value = 2.6
$ws.Cells.Item(202,14).Value = 2.25
$ws.Cells.Item(202,15).Value = 3.3
$ws.Cells.Item(202,16).Value = 2.8
$ws.Cells.Item(202,17).Value = -0.25
$ws.Cells.Item(202,18).Value = 2
$ws.Cells.Item(202,19).Value = 1.8
$ws.Cells.Item(202,20).Value = 2.5
$ws.Cells.Item(202,21).Value = 1.85
$ws.Cells.Item(202,22).Value = 1.95
$ws.Cells.Item(202,23).Value = 0
$ws.Cells.Item(202,24).Value = 0
$ws.Cells.Item(202,25).Value = 0
$ws.Cells.Item(202,26).Value = 0
$ws.Cells.Item(202,27).Value = 0

# Row 203
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(203,1).PasteSpecial(-4122)
$ws.Cells.Item(203,1).Value = 201
$ws.Cells.Item(203,2).Value = 7623987
$ws.Cells.Item(203,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(203,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(203,5).PasteSpecial(-4122)
$ws.Cells.Item(203,5).Value = 45354.79166666666
$ws.Cells.Item(203,6).Value = 'Deportivo Saprissa'
$ws.Cells.Item(203,7).Value = 'AD San Carlos'
$ws.Cells.Item(203,11).Value = 1.571
$ws.Cells.Item(203,12).Value = 3.75
$ws.Cells.Item(203,13).Value = 5
$ws.Cells.Item(203,14).Value = 1.533
$ws.Cells.Item(203,15).Value = 3.8
$ws.Cells.Item(203,16).Value = 5.25
$ws.Cells.Item(203,17).Value = -1
$ws.Cells.Item(203,18).Value = 1.975
$ws.Cells.Item(203,19).Value = 1.825
$ws.Cells.Item(203,20).Value = 2.5
$ws.Cells.Item(203,21).Value = 1.825
$ws.Cells.Item(203,22).Value = 1.975
$ws.Cells.Item(203,23).Value = 0
$ws.Cells.Item(203,24).Value = 0
$ws.Cells.Item(203,25).Value = 0
$ws.Cells.Item(203,26).Value = 0
$ws.Cells.Item(203,27).Value = 0

# Row 204
$ws.Cells.Item(198,1).Copy()
$ws.Cells.Item(204,1).PasteSpecial(-4122)
$ws.Cells.Item(204,1).Value = 202
$ws.Cells.Item(204,2).Value = 7623988
$ws.Cells.Item(204,3).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(204,4).Value = 'Costa Rica Primera Division'
$ws.Cells.Item(198,5).Copy()
$ws.Cells.Item(204,5).PasteSpecial(-4122)
$ws.Cells.Item(204,5).Value = 45354.875
$ws.Cells.Item(204,6).Value = 'Sporting San Jose'
$ws.Cells.Item(204,7).Value = 'Cartagines'
$ws.Cells.Item(204,11).Value = 2.45
$ws.Cells.Item(204,12).Value = 3.25
$ws.Cells.Item(204,13).Value = 2.6
$ws.Cells.Item(204,14).Value = 2.45
$ws.Cells.Item(204,15).Value = 3.25
$ws.Cells.Item(204,16).Value = 2.6
$ws.Cells.Item(204,17).Value = 0
$ws.Cells.Item(204,18).Value = 1.825
$ws.Cells.Item(204,19).Value = 1.975
$ws.Cells.Item(204,20).Value = 2.5
$ws.Cells.Item(204,21).Value = 1.875
$ws.Cells.Item(204,22).Value = 1.925
$ws.Cells.Item(204,23).Value = 0
$ws.Cells.Item(204,24).Value = 0
$ws.Cells.Item(204,25).Value = 0
$ws.Cells.Item(204,26).Value = 0
$ws.Cells.Item(204,27).Value = 0
